$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- A1: documentation blurb about automatic MATERIAL mapping ---
$ws.Range("A1").Value = "This column will automatically map to the MATERIAL data column because the name matches (case sensitive) and there is no explicit mapping provided in the sample Powermax105 XML transform."
$ws.Range("A1").WrapText = $true
$ws.Range("A1").Font.Name = "Tahoma"
$ws.Range("A1").Font.Size = 8

# --- G1: documentation blurb about BASE_FEEDRATE mapping, now with rich-text emphasis ---
$g1Text = "This column is mapped to the BASE_FEEDRATE data column by the header attribute in the XML transform."
$ws.Range("G1").Value = $g1Text
$ws.Range("G1").WrapText = $true

# " XML" segment (characters 86-89, 1-based)
$runXml = $ws.Range("G1").Characters(86, 4)
$runXml.Font.Name = "Tahoma"
$runXml.Font.Size = 8
$runXml.Font.ColorIndex = -4105

# " " segment between XML and transform - bold & red
$runBoldRed = $ws.Range("G1").Characters(90, 1)
$runBoldRed.Font.Name = "Tahoma"
$runBoldRed.Font.Size = 8
$runBoldRed.Font.Bold = $true
$runBoldRed.Font.Color = 255

# "transform." segment
$runTransform = $ws.Range("G1").Characters(91, 10)
$runTransform.Font.Name = "Tahoma"
$runTransform.Font.Size = 8

# --- J1: documentation blurb about column width ---
$ws.Range("J1").Value = "This column's width will be adjusted by the width attribute in the XML transform."
$ws.Range("J1").WrapText = $true
$ws.Range("J1").Font.Name = "Tahoma"
$ws.Range("J1").Font.Size = 8

# --- U1: unchanged text, left as-is ---
$ws.Range("U1").Value = "Unmapped data columns and unmapped custom columns will appear here and to the right."
$ws.Range("U1").WrapText = $true

# Row 1 auto-fits a bit taller with the revised wording
$ws.Rows.Item(1).RowHeight = 61.2

# Move the saved cursor position off F2, matching the refreshed view
[void]$ws.Range("A1").Select()
